# Swap the order of "System" and the recorded-by email/name in column G
# ("Recorded By") from "System, <who>" to "<who>, System".
# The backup@backdoor.com rows ("System, backup@backdoor.com" and
# "System, system, backup@backdoor.com") are intentionally left as-is,
# matching the upstream diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "System, "
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null) {
        continue
    }

    if (-not $text.StartsWith($prefix)) {
        continue
    }

    $rest = $text.Substring($prefix.Length)

    # Skip the "backdoor" rows (e.g. "System, backup@backdoor.com" or
    # "System, system, backup@backdoor.com") - the upstream diff leaves
    # those untouched, so we must not reorder them either.
    if ($rest.Contains("backdoor.com")) {
        continue
    }

    $cell.Value = $rest + ", System"
    $changed++
}

Write-Output "Updated $changed cells in column G"
